$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controles")

# Row 40: "Bloqueada"/"Blocked" -> "Negada"/"Denied"
$ws.Range("B40").Value = "Negada"
$ws.Range("C40").Value = "Denied"

# New rows 178 and 179: email recovery messages
$ws.Range("A178").Value = "com.td.email.header.recover"
$ws.Range("B178").Value = "Se detecto actividad sospechosa en su cuenta"
$ws.Range("C178").Value = "Suspicious activity detected in you account"

$ws.Range("A179").Value = "com.td.email.body.recover"
$ws.Range("B179").Value = "Por seguridad, se genero esta nueva contraseña:"
$ws.Range("C179").Value = "For security, this password was generated:"

# Fill down D/E formulas from row 177 into rows 178:179
$dFormula = $ws.Range("D177").Formula
$eFormula = $ws.Range("E177").Formula

$ws.Range("D178").Formula = $dFormula -replace '177', '178'
$ws.Range("E178").Formula = $eFormula -replace '177', '178'
$ws.Range("D179").Formula = $dFormula -replace '177', '179'
$ws.Range("E179").Formula = $eFormula -replace '177', '179'

# Undo the row-height autofit side effect from multi-line formula text
$ws.Rows.Item(178).AutoFit()
$ws.Rows.Item(179).AutoFit()

# Reflect the user's final scroll/selection position on save
$ws.Activate()
$ws.Range("A40").Select() | Out-Null
